$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lua")

# Insert a new row above row 8. This shifts the "category" labels in
# column A (and their per-row styling), together with the C:H "count"
# block further down, one row down - exactly like Excel's native
# "Insert Sheet Rows" command - and leaves a freshly (and correctly)
# formatted blank row 8 behind (styles s=1/7/6 for A8/B8/C8).
$ws.Rows.Item(8).Insert()

# Fill in the newly inserted row with the "country-codes" entry.
$ws.Cells.Item(8, 1).Value2 = "country-codes"
$ws.Cells.Item(8, 2).Value2 = 7

# The "from/to/type" numbering in column B is keyed off the row
# position, not carried along with the row insert, so bump B9:B13 by
# one and extend the sequence onto the row (14) that used to hold the
# "include-float-in-format" count block before it shifted down.
$ws.Cells.Item(9, 2).Value2 = 8
$ws.Cells.Item(10, 2).Value2 = 9
$ws.Cells.Item(11, 2).Value2 = 10
$ws.Cells.Item(12, 2).Value2 = 11
$ws.Cells.Item(13, 2).Value2 = 12
$ws.Cells.Item(14, 2).Value2 = 13

# Update the selection to match the saved state of the workbook.
$null = $ws.Range("B15").Select()
